$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 5 with the next id and sample data, matching the style of the
# existing id column (A2:A4 uses the same style as row 4's id cell).
$ws.Range("A4").Copy($ws.Range("A5"))
$ws.Range("A5").Value = 3

$ws.Range("B5").Value = "ss"
$ws.Range("C5").Value = "ss"
$ws.Range("D5").Value = "sss"
$ws.Range("E5").Value = "sss"
$ws.Range("F5").Value = "sss"
$ws.Range("G5").Value = "ss"
